$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates derived from the authoritative diff.
# Column D (Price) values are forced to Text format to preserve
# the exact string representation (e.g. "1.00", "247.45") since
# Excel would otherwise auto-convert numeric-looking text to numbers.
$updates = [ordered]@{
    "D2" = "35.430.44"
    "E2" = "  +0.43%  "
    "D3" = "1.894.07"
    "E3" = "  -0.79%  "
    "E4" = "  -0.73%  "
    "D5" = "247.45"
    "E5" = "  -3.32%  "
    "E6" = "  -4.65%  "
    "E7" = "  -0.83%  "
    "D8" = "43.95"
    "E8" = "  +7.88%  "
    "D9" = "0.354"
    "E9" = "  -4.86%  "
    "D10" = "0.0740"
    "E10" = "  -2.46%  "
    "D11" = "0.0970"
    "E11" = "  -1.84%  "
    "D12" = "13.08"
    "E12" = "  +1.82%  "
    "D13" = "2.168.61"
    "E13" = "  -0.86%  "
    "E14" = "  -0.75%  "
    "E15" = "  -0.90%  "
    "D16" = "1.901.13"
    "E16" = "  -0.05%  "
    "D17" = "35.419.09"
    "E17" = "  +0.40%  "
    "D18" = "73.68"
    "E18" = "  -1.23%  "
    "D19" = "0.0₃0822"
    "E19" = "  -3.31%  "
    "D20" = "247.64"
    "E20" = "  +1.55%  "
    "D21" = "12.82"
    "E21" = "  -1.77%  "
    "D22" = "4.95"
    "E22" = "  -3.10%  "
    "D23" = "1.00"
    "E23" = "  -0.76%  "
    "D24" = "2.55"
    "E24" = "  +5.58%  "
    "E25" = "  -10.19%  "
    "D26" = "165.49"
    "E26" = "  -0.63%  "
    "D27" = "8.45"
    "E27" = "  -2.79%  "
    "D28" = "18.38"
    "E28" = "  -1.96%  "
    "E29" = "  -3.69%  "
    "D30" = "4.128.41"
    "E30" = "  -0.02%  "
    "E31" = "  +7.87%  "
    "E32" = "  -3.07%  "
    "D33" = "0.0579"
    "E33" = "  -1.60%  "
    "E34" = "  -0.17%  "
    "E35" = "  -0.81%  "
    "D36" = "0.858"
    "E36" = "  -5.76%  "
    "E37" = "  -1.47%  "
    "D38" = "1.59"
    "E38" = "  -20.73%  "
    "D39" = "0.0694"
    "E39" = "  +6.45%  "
    "B40" = "Aave"
    "C40" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D40" = "97.64"
    "E40" = "  +1.13%  "
    "B41" = "InjectiveProtocol"
    "C41" = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
    "D41" = "17.06"
    "E41" = "  -1.39%  "
    "E42" = "  -2.54%  "
    "E43" = "  -2.76%  "
    "B44" = "Maker"
    "C44" = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
    "D44" = "1.293.90"
    "E44" = "  -3.24%  "
    "B45" = "RenderToken"
    "C45" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D45" = "2.36"
    "E45" = "  -2.80%  "
    "D46" = "0.0798"
    "E46" = "  +6.51%  "
    "E47" = "  -1.00%  "
    "D48" = "2.75"
    "E48" = "  -0.54%  "
    "D49" = "12.12"
    "E49" = "  +3.92%  "
    "D50" = "6.37"
    "E50" = "  -5.26%  "
    "D51" = "43.17"
    "E51" = "  -3.99%  "
}

foreach ($addr in $updates.Keys) {
    $col = $addr -replace '[0-9]+$', ''
    $cell = $ws.Range($addr)
    if ($col -eq "D") {
        # Force text storage so numeric-looking strings keep their
        # exact formatting (leading/trailing zeros, no auto-number coercion).
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$addr]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$addr]
    }
}

